$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 8 ("1999-2001" / 埼玉県衛生研究所 ...)
# so that row shifts down to row 9, and the new row 8 can hold the
# "2004-2011" entry that was added to the source Markdown table.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "2004-2011"
$ws.Range("B8").Value = "**埼玉県衛生研究所** <br> [市販鶏肉のカンピロバクター及びサルモネラ汚染状況と分離株の薬剤感受性](https://jvma-vet.jp/mag/06706/d2.pdf) <br>（日獣会誌, 67, 442~448, 2014）"
$ws.Range("C8").Value = "済"
